# Quality_Assurance_KPI_CHC.xlsx — "update sample templates downloads for
# HP Kayakalp-program"
#
# The authored diff shows the KPI value column (C) being wiped of its
# sample/placeholder numbers (20.6, 2, 3, 4, ... 35) for every data row,
# while the labels in columns A/B and the section headers stay untouched.
# It also shows the saved view state (scrolled position / active cell)
# moving further down the sheet, which happens naturally when someone
# scrolls down and clicks a cell before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the sample numeric values out of the KPI value column (C4:C43).
# The section-header rows (11, 17, 26, 33) and the sub-header rows
# (37, 38) in that range are already blank, so clearing the whole
# contiguous block reproduces exactly the set of cells touched by the
# diff without disturbing any of the surrounding labels/styles.
$ws.Range("C4:C43").ClearContents()

# Reflect the saved scroll/selection state recorded in the diff
# (selection moved from D36 to D48, further down the sheet).
[void]$ws.Range("D48").Select()
